# "Generate Report for Handback"
#
# The c98032de-8112-41da-ad19-ff7e82fdb0f9.md file (rows 3 & 4 in every
# sheet) has now been handed back. Update the Overview sheet and the two
# per-locale sheets (zh-cn, de-de) to reflect this: flip the Status from
# "Ready for handoff" to "Handed back: in sync with en-US", and populate
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns that were previously empty placeholders.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack
$overview.Range("B4").Value = $statusHandedBack
$overview.Range("C4").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("C4").Value = $statusHandedBack

$zhcn.Range("F3").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.md"
$zhcn.Range("G3").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-03-19 16:51:00"

$zhcn.Range("F4").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.md"
$zhcn.Range("G4").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf"
$zhcn.Range("H4").Value = "2016-03-19 16:51:00"

$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/137d1158a49d00463cd556821ca460a443ffc477/e2e/c98032de-8112-41da-ad19-ff7e82fdb0f9.md", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a32a0651cc6006f4ab7dcdcefa3185ebbdec49d2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/137d1158a49d00463cd556821ca460a443ffc477/e2e/c98032de-8112-41da-ad19-ff7e82fdb0f9.md", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a32a0651cc6006f4ab7dcdcefa3185ebbdec49d2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("C4").Value = $statusHandedBack

$dede.Range("F3").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.md"
$dede.Range("G3").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf"
$dede.Range("H3").Value = "2016-03-19 16:51:14"

$dede.Range("F4").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.md"
$dede.Range("G4").Value = "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf"
$dede.Range("H4").Value = "2016-03-19 16:51:14"

$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/137d1158a49d00463cd556821ca460a443ffc477/e2e/c98032de-8112-41da-ad19-ff7e82fdb0f9.md", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e25e69183f7a6550357526cf0bc4a08a2438cc27/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/137d1158a49d00463cd556821ca460a443ffc477/e2e/c98032de-8112-41da-ad19-ff7e82fdb0f9.md", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.md")
$dede.Hyperlinks.Add($dede.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e25e69183f7a6550357526cf0bc4a08a2438cc27/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf", [Type]::Missing, [Type]::Missing, "c98032de-8112-41da-ad19-ff7e82fdb0f9.dfaa3266185d154a54583c07bd5c0ae8ba75c240.de-de.xlf")
